$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Job 3's JSON payload, updated so "active" is now 0 (was 1)
$job3Json = '[{"job_id":3,"customer":"Alex","site":"Brighton Avenue","cs_number":"666","dwg_number":"0","start_date":19180,"required_by":19208,"hours_mill":25,"hours_program":10,"hours_cnc":20,"hours_veneer":0,"hours_bench":10,"hours_spray":0,"hours_dispatch":1,"active":0,"created_by":"a117644@r02.xlgs.local","created_on":"2022-07-07 17:37:37","updated_by":"a117644@r02.xlgs.local","updated_on":"2022-07-07 17:37:37","start_date_mill":19180,"days_mill":4,"end_date_mill":19186,"start_date_program":19186,"days_program":2,"end_date_program":19188,"start_date_cnc":19188,"days_cnc":3,"end_date_cnc":19193,"start_date_veneer":19193,"days_veneer":0,"end_date_veneer":19193,"start_date_bench":19193,"days_bench":2,"end_date_bench":19195,"start_date_spray":19195,"days_spray":0,"end_date_spray":19195,"start_date_dispatch":19195,"days_dispatch":1,"end_date_dispatch":19198,"days_overrun":0,"start_date_recommended":19180}]'

# Update existing row 4 (job_id 3) data string to reflect the new "active" value
$ws.Range("B4").Value = $job3Json

# Add a new row 5 duplicating job_id 3's record with the updated data,
# the start of incorporating data sourced from Google Sheets
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = $job3Json
